$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles, row height, number formats) from the last existing row
$ws.Range("A36:C36").Copy()
$ws.Range("A37:C37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in today's journal entry on the new row (row 37)
$ws.Range("A37").Value = (Get-Date -Year 2018 -Month 3 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B37").Value = "J'ai avancé la partie du code pour ajouter un article dans la base de données. Je dois encore ajouter des paramètres tels que le prix, la taille et la couleur dans la BD. Une fois que tout fonctionne je vais regarder si je peux optimiser mes requêtes."
$ws.Range("C37").Value = "3 périodes"

$ws.Rows.Item(37).RowHeight = 45

$ws.Range("C38").Select()
